$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRIM")
$ws.Range("M5").Value = "test"
